$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("統計")

$ws.Cells.Item(8, 1).Value = "2025-08-28T01:45:37.372586"
$ws.Cells.Item(8, 2).Value = 13
$ws.Cells.Item(8, 3).Value = "全案件リスト"
$ws.Cells.Item(8, 4).Value = 69.2
$ws.Cells.Item(8, 5).Value = 4
$ws.Cells.Item(8, 6).Value = 6
$ws.Cells.Item(8, 7).Value = 13
